$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained an extra row: the "Accuracy over PyType" summary that used
# to live in row 3 (columns E/F) moves down to row 4, and row 3 instead gains
# a new "Scalpel Accuracy:" summary in columns C/D.

# 1. Carry the existing E3/F3 values down to the new row 4, matching the
#    formatting (white-fill style) already used across row 3.
$ws.Range("A4:F4").Interior.ColorIndex = 2
$ws.Range("A4:D4").Value2 = ""
$ws.Range("E4").Value2 = $ws.Range("E3").Value2
$ws.Range("F4").Value2 = $ws.Range("F3").Value2

# 2. Clear out E3/F3 now that their content has moved to row 4.
$ws.Range("E3").Value2 = ""
$ws.Range("F3").Value2 = ""

# 3. Populate the new "Scalpel Accuracy:" summary cells in row 3.
$ws.Range("C3").Value2 = "Scalpel Accuracy:"
$ws.Range("D3").Value2 = 100
